# Applies the cryptos.xlsx price/volume/listing update described by the diff.
# Uses NumberFormat "@" + Style reset so numeric-looking price strings (e.g. "516.60")
# are stored as text (matching the original inlineStr/shared-string text cells)
# instead of being auto-converted to numbers by Excel, while keeping the cell style
# identical to the original (no explicit style index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.331.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("E10").Value = "  -0.64%  "

$ws.Range("E11").Value = "  -1.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.599.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.426.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.081.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.501"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0911"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.66%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.21%  "

$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0673"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.112.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.658"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.269.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "

$ws.Range("E45").Value = "  +7.63%  "

$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.933"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0872"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "249.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.12%  "
